# Applies the scheduled market-data refresh to the Chocobo Profits workbook.
# For every Leve row on every job sheet (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR) this
# rewrites the live-market derived columns (H:N -- current average prices,
# leve NQ/HQ prices, and NQ/HQ profit) with freshly polled values. A couple of
# rows gain or lose their profit columns entirely when the market data for that
# item newly has/lacks a price.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 1873.4286
$ws.Range("I96").Value = 622.8
$ws.Range("J96").Value = 5000
$ws.Range("K96").Value = 1868.4
$ws.Range("L96").Value = 15000
$ws.Range("M96").Value = -495.3999999999999
$ws.Range("N96").Value = -17746
$ws.Range("H112").Value = 19609436
$ws.Range("I112").Value = 333334080
$ws.Range("J112").Value = 1646
$ws.Range("K112").Value = 1000002240
$ws.Range("L112").Value = 4938
$ws.Range("M112").Value = -1000001132
$ws.Range("N112").Value = -7154
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("H129").Value = 649.7778
$ws.Range("I129").Value = 371.83334
$ws.Range("J129").Value = 1205.6666
$ws.Range("K129").Value = 1115.50002
$ws.Range("L129").Value = 3616.9998
$ws.Range("M129").Value = 3884.49998
$ws.Range("N129").Value = -13616.9998
$ws.Range("H132").Value = 27033226
$ws.Range("I132").Value = 33339316
$ws.Range("J132").Value = 7130.2856
$ws.Range("K132").Value = 100017948
$ws.Range("L132").Value = 21390.8568
$ws.Range("M132").Value = -100015418
$ws.Range("N132").Value = -26450.8568
$ws.Range("H135").Value = 313.56412
$ws.Range("I135").Value = 313.56412
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 2822.07708
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -287.07708
$ws.Range("H137").Value = 2964.4844
$ws.Range("I137").Value = 2503.054
$ws.Range("J137").Value = 3596.8147
$ws.Range("K137").Value = 7509.162
$ws.Range("L137").Value = 10790.4441
$ws.Range("M137").Value = -4959.162
$ws.Range("N137").Value = -15890.4441
$ws.Range("H138").Value = 2923.82
$ws.Range("I138").Value = 487.39285
$ws.Range("J138").Value = 3871.3193
$ws.Range("K138").Value = 1462.17855
$ws.Range("L138").Value = 11613.9579
$ws.Range("M138").Value = 3677.82145
$ws.Range("N138").Value = -21893.9579
$ws.Range("N120").Value = $null
$ws.Range("N135").Value = $null

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3859.023
$ws.Range("I32").Value = 3932.7742
$ws.Range("K32").Value = 3932.7742
$ws.Range("M32").Value = -3645.7742
$ws.Range("H61").Value = 1657.3448
$ws.Range("I61").Value = 1249.7059
$ws.Range("J61").Value = 2234.8333
$ws.Range("K61").Value = 1249.7059
$ws.Range("L61").Value = 2234.8333
$ws.Range("M61").Value = -1037.7059
$ws.Range("N61").Value = -2658.8333
$ws.Range("H74").Value = 2564.6223
$ws.Range("I74").Value = 2609.7778
$ws.Range("J74").Value = 2384
$ws.Range("K74").Value = 2609.7778
$ws.Range("L74").Value = 2384
$ws.Range("M74").Value = -1735.7778
$ws.Range("N74").Value = -4132
$ws.Range("H77").Value = 2564.6223
$ws.Range("I77").Value = 2609.7778
$ws.Range("J77").Value = 2384
$ws.Range("K77").Value = 13048.889
$ws.Range("L77").Value = 11920
$ws.Range("M77").Value = -8680.888999999999
$ws.Range("N77").Value = -20656
$ws.Range("H109").Value = 26050
$ws.Range("J109").Value = 26050
$ws.Range("L109").Value = 26050
$ws.Range("N109").Value = -28824
$ws.Range("H132").Value = 2352.8
$ws.Range("I132").Value = 1684.8334
$ws.Range("J132").Value = 4070.4285
$ws.Range("K132").Value = 5054.5002
$ws.Range("L132").Value = 12211.2855
$ws.Range("M132").Value = -2524.5002
$ws.Range("N132").Value = -17271.2855
$ws.Range("H136").Value = 1657.3448
$ws.Range("I136").Value = 1249.7059
$ws.Range("J136").Value = 2234.8333
$ws.Range("K136").Value = 3749.1177
$ws.Range("L136").Value = 6704.499899999999
$ws.Range("M136").Value = -1199.1177
$ws.Range("N136").Value = -11804.4999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2005.0735
$ws.Range("I134").Value = 1211.2858
$ws.Range("J134").Value = 4052.2104
$ws.Range("K134").Value = 3633.8574
$ws.Range("L134").Value = 12156.6312
$ws.Range("M134").Value = -1098.8574
$ws.Range("N134").Value = -17226.6312
$ws.Range("H137").Value = 46583.75
$ws.Range("J137").Value = 53778.332
$ws.Range("L137").Value = 53778.332
$ws.Range("N137").Value = -63978.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 6946535.5
$ws.Range("I16").Value = 12347362
$ws.Range("J16").Value = 2615.1428
$ws.Range("K16").Value = 12347362
$ws.Range("L16").Value = 2615.1428
$ws.Range("M16").Value = -12347075
$ws.Range("N16").Value = -3189.1428
$ws.Range("H31").Value = 8774282
$ws.Range("I31").Value = 1361.711
$ws.Range("J31").Value = 41672736
$ws.Range("K31").Value = 1361.711
$ws.Range("L31").Value = 41672736
$ws.Range("M31").Value = -1066.711
$ws.Range("N31").Value = -41673326
$ws.Range("H34").Value = 8774282
$ws.Range("I34").Value = 1361.711
$ws.Range("J34").Value = 41672736
$ws.Range("K34").Value = 1361.711
$ws.Range("L34").Value = 41672736
$ws.Range("M34").Value = -1159.711
$ws.Range("N34").Value = -41673140
$ws.Range("H58").Value = 1366.84
$ws.Range("I58").Value = 1470.3158
$ws.Range("J58").Value = 1039.1666
$ws.Range("K58").Value = 1470.3158
$ws.Range("L58").Value = 1039.1666
$ws.Range("M58").Value = -1267.3158
$ws.Range("N58").Value = -1445.1666
$ws.Range("H104").Value = 31935
$ws.Range("I104").Value = 21000
$ws.Range("J104").Value = 33028.5
$ws.Range("K104").Value = 21000
$ws.Range("L104").Value = 33028.5
$ws.Range("M104").Value = -18379
$ws.Range("N104").Value = -38270.5
$ws.Range("H113").Value = 6946535.5
$ws.Range("I113").Value = 12347362
$ws.Range("J113").Value = 2615.1428
$ws.Range("K113").Value = 12347362
$ws.Range("L113").Value = 2615.1428
$ws.Range("M113").Value = -12345192
$ws.Range("N113").Value = -6955.1428
$ws.Range("H134").Value = 2856.5908
$ws.Range("I134").Value = 3221.875
$ws.Range("J134").Value = 2294.6155
$ws.Range("K134").Value = 9665.625
$ws.Range("L134").Value = 6883.8465
$ws.Range("M134").Value = -7130.625
$ws.Range("N134").Value = -11953.8465
$ws.Range("H136").Value = 1366.84
$ws.Range("I136").Value = 1470.3158
$ws.Range("J136").Value = 1039.1666
$ws.Range("K136").Value = 4410.9474
$ws.Range("L136").Value = 3117.4998
$ws.Range("M136").Value = -1860.9474
$ws.Range("N136").Value = -8217.4998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 583
$ws.Range("I113").Value = 490.87234
$ws.Range("J113").Value = 892.2857
$ws.Range("K113").Value = 1472.61702
$ws.Range("L113").Value = 2676.8571
$ws.Range("M113").Value = 697.3829799999999
$ws.Range("N113").Value = -7016.8571

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3286.1562
$ws.Range("I132").Value = 2209.2778
$ws.Range("J132").Value = 4670.7144
$ws.Range("K132").Value = 6627.8334
$ws.Range("L132").Value = 14012.1432
$ws.Range("M132").Value = -4097.8334
$ws.Range("N132").Value = -19072.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4912.1875
$ws.Range("I7").Value = 3465.5557
$ws.Range("J7").Value = 6772.143
$ws.Range("K7").Value = 3465.5557
$ws.Range("L7").Value = 6772.143
$ws.Range("M7").Value = -3353.5557
$ws.Range("N7").Value = -6996.143
$ws.Range("H69").Value = 39999
$ws.Range("J69").Value = 39999
$ws.Range("L69").Value = 39999
$ws.Range("N69").Value = -41621
$ws.Range("H72").Value = 39999
$ws.Range("J72").Value = 39999
$ws.Range("L72").Value = 119997
$ws.Range("N72").Value = -128109
$ws.Range("H122").Value = 6600.6665
$ws.Range("I122").Value = 3201
$ws.Range("J122").Value = 8300.5
$ws.Range("K122").Value = 9603
$ws.Range("L122").Value = 24901.5
$ws.Range("M122").Value = -7153
$ws.Range("N122").Value = -29801.5
$ws.Range("H126").Value = 4912.1875
$ws.Range("I126").Value = 3465.5557
$ws.Range("J126").Value = 6772.143
$ws.Range("K126").Value = 10396.6671
$ws.Range("L126").Value = 20316.429
$ws.Range("M126").Value = -7926.667099999999
$ws.Range("N126").Value = -25256.429
$ws.Range("H132").Value = 3109.5715
$ws.Range("I132").Value = 1015.1875
$ws.Range("J132").Value = 6576.1377
$ws.Range("K132").Value = 3045.5625
$ws.Range("L132").Value = 19728.4131
$ws.Range("M132").Value = -515.5625
$ws.Range("N132").Value = -24788.4131
$ws.Range("H136").Value = 2425.9812
$ws.Range("I136").Value = 1538
$ws.Range("J136").Value = 4306.4116
$ws.Range("K136").Value = 4614
$ws.Range("L136").Value = 12919.2348
$ws.Range("M136").Value = -2064
$ws.Range("N136").Value = -18019.2348

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 39583.668
$ws.Range("I80").Value = 30600.5
$ws.Range("J80").Value = 44075.25
$ws.Range("K80").Value = 30600.5
$ws.Range("L80").Value = 44075.25
$ws.Range("M80").Value = -29602.5
$ws.Range("N80").Value = -46071.25
$ws.Range("H83").Value = 39583.668
$ws.Range("I83").Value = 30600.5
$ws.Range("J83").Value = 44075.25
$ws.Range("K83").Value = 91801.5
$ws.Range("L83").Value = 132225.75
$ws.Range("M83").Value = -86809.5
$ws.Range("N83").Value = -142209.75
$ws.Range("H122").Value = 2063.15
$ws.Range("I122").Value = 1286.3226
$ws.Range("J122").Value = 4738.8887
$ws.Range("K122").Value = 3858.9678
$ws.Range("L122").Value = 14216.6661
$ws.Range("M122").Value = -1408.9678
$ws.Range("N122").Value = -19116.6661
$ws.Range("H132").Value = 4169074.2
$ws.Range("I132").Value = 2402.8
$ws.Range("J132").Value = 13335752
$ws.Range("K132").Value = 7208.400000000001
$ws.Range("L132").Value = 40007256
$ws.Range("M132").Value = -4678.400000000001
$ws.Range("N132").Value = -40012316
$ws.Range("H136").Value = 3231.3857
$ws.Range("I136").Value = 3405.5945
$ws.Range("J136").Value = 3036.0605
$ws.Range("K136").Value = 10216.7835
$ws.Range("L136").Value = 9108.181500000001
$ws.Range("M136").Value = -7666.783500000001
$ws.Range("N136").Value = -14208.1815
